$wb = $excel.ActiveWorkbook

# Banner_Text sheet holds the message table that needs updated wording.
$banner = $wb.Worksheets.Item("Banner_Text")

$banner.Range("B3").Value = "శ్రీ వేదాంత దేశికులచే క్రి.శ. 1359వ సంవత్సరంలో స్థాపితమైనది"
$banner.Range("B3").Style = "Normal"

$banner.Range("B4").Value = "శ్రీ బ్రహ్మతంత్ర స్వతంత్ర పరకాల స్వామి మఠ ఆచార్యులు"
$banner.Range("B4").Style = "Normal"

$banner.Range("B6").Value = "శ్రీ పరకాల స్వామి మఠ – శ్రీ వేదాంత దేశికుల అవిచ్ఛిన్న పరంపర"
$banner.Range("B6").Style = "Normal"

$banner.Range("B2").Value = "శ్రీ బ్రహ్మతంత్ర స్వతంత్ర పరకాల స్వామి మఠ గురుపరంపర"
$banner.Range("B2").Style = "Normal"

# Column C on Banner_Text widened (and best-fit) in the author's final pass.
$banner.Columns.Item(3).ColumnWidth = 52.7265625

# Selection/active-cell moves to B6, and Banner_Text becomes the active tab.
$banner.Activate()
$banner.Range("B6").Select()
